$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text entry (not numeric/percent auto-conversion), then restore
# the default "Normal" cell style so no stray style index is left behind.

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "{'criterion': 'gini', 'max_depth': 15, 'min_samples_leaf': 2, 'min_samples_split': 15}"
$ws.Range("C2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "62.58%"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "63.81%"
$ws.Range("F2").Style = "Normal"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "62.58%"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "62.40%"
$ws.Range("H2").Style = "Normal"
